$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header row: Date, Calories, Fat (g), Fiber (g), Carbs (g), Sodium (mg),
# Protein (g), Water (ml), Energy From Carbs, Energy From Fat,
# Energy From Protein
# ---------------------------------------------------------------------
$headers = @("Date", "Calories", "Fat (g)", "Fiber (g)", "Carbs (g)", "Sodium (mg)", "Protein (g)", "Water (ml)", "Energy From Carbs", "Energy From Fat", "Energy From Protein")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# ---------------------------------------------------------------------
# Data rows 2-6: Date, Calories, Fat, Fiber, Carbs, Sodium, Protein,
# Water, % energy from carbs, % energy from fat, % energy from protein
# ---------------------------------------------------------------------
$dates = @("10/27/2023", "10/28/2023", "10/29/2023", "10/30/2023", "10/31/2023")

$numbers = @(
    @(1323, 46.1, 0, 116.2, 1477.6, 105.7, 0),
    @(1498, 51.2, 0, 103.5, 1863.2, 152.4, 0),
    @(1492, 46.5, 0, 195,   475.7,  73.6,  0),
    @(1327, 64.3, 0, 109.4, 1137.9, 75.6,  0),
    @(1304, 37.6, 0, 129.6, 899.4,  112.2, 0)
)

$pcts = @(
    @("36%", "32%", "32%"),
    @("28%", "31%", "41%"),
    @("52%", "28%", "20%"),
    @("33%", "44%", "23%"),
    @("40%", "26%", "34%")
)

for ($i = 0; $i -lt 5; $i++) {
    $row = $i + 2

    # Date goes through a TEXT() formula so it lands verbatim as text
    # instead of Excel's automatic "looks like a date -> date serial"
    # conversion; converted to a static value below.
    $ws.Cells.Item($row, 1).Formula = '=TEXT("' + $dates[$i] + '","@")'

    # Plain numeric columns B-H.
    $nums = $numbers[$i]
    for ($c = 0; $c -lt $nums.Length; $c++) {
        $ws.Cells.Item($row, 2 + $c).Value = $nums[$c]
    }

    # Percentage columns I-K, same TEXT() trick so "36%" etc. is kept as
    # literal text rather than becoming the number 0.36 with a percent
    # format.
    $rowPcts = $pcts[$i]
    $ws.Cells.Item($row, 9).Formula  = '=TEXT("' + $rowPcts[0] + '","@")'
    $ws.Cells.Item($row, 10).Formula = '=TEXT("' + $rowPcts[1] + '","@")'
    $ws.Cells.Item($row, 11).Formula = '=TEXT("' + $rowPcts[2] + '","@")'
}

# Freeze the TEXT() formulas into plain text values (Copy + PasteSpecial
# values-only) so the saved cells are plain shared strings with no
# formula and no left-over number-format style change.
$dateRange = $ws.Range("A2:A6")
$dateRange.Copy()
$dateRange.PasteSpecial(-4163)

$pctRange = $ws.Range("I2:K6")
$pctRange.Copy()
$pctRange.PasteSpecial(-4163)
